# Loan RBI, Variable Instalments
#
# On the "Repayment schedule" sheet, insert a new (blank) column before
# column N, pushing the old N/O/P ("In Advance" / "Paid Date" / "Outstanding")
# columns right by one into O/P/Q. Give the freshly inserted column the same
# width as the neighbouring "Due" column (11 characters, no bestFit).
#
# Also switch the active/selected sheet from "Edit Repayment Schedule" to
# "Repayment schedule", with the selection on the latter moved to I15.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

$ws.Columns("N").Insert()
$ws.Columns("N").ColumnWidth = 10.17

[void]$ws.Activate()
[void]$ws.Range("I15").Select()
